$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-13 (columns A:T), combining the existing M1->M2 cluster rows
# (recomputed values) with six newly added M2-sending-cluster rows, per the
# updated Natmi analysis.
$data = @(
        @("M1", "Ifnk", "Ifnar2", "ECs", 1, 0.3333333333333333, 0.1379033333333333, 0.41371, 0.1151434304008604, 0.1151434304008603, 2, 1, 29.4475095, 58.895019, 0.1188040616343076, 0.08432571340736816, 4.060909718415, 24.36545831049, 0.01367950720212942, 0.009709551912724193),
        @("M1", "Ifnk", "Ifnar2", "FAPs", 1, 0.3333333333333333, 0.1379033333333333, 0.41371, 0.1151434304008604, 0.1151434304008603, 3, 1, 24.18563533333333, 72.556906, 0.09757537257289008, 0.1038867626663191, 3.335279731251112, 30.01751758126, 0.01123516312068459, 0.01196187822664001),
        @("M1", "Ifnk", "Ifnar2", "M1", 1, 0.3333333333333333, 0.1379033333333333, 0.41371, 0.1151434304008604, 0.1151434304008603, 3, 1, 73.143293, 219.429879, 0.2950918579280267, 0.3141790494976744, 10.08670391567667, 90.78033524109, 0.03397788880519632, 0.03617565351924393),
        @("M1", "Ifnk", "Ifnar2", "M2", 1, 0.3333333333333333, 0.1379033333333333, 0.41371, 0.1151434304008604, 0.1151434304008603, 3, 1, 82.44778433333333, 247.343353, 0.3326302229921831, 0.3541454787253825, 11.36982428551445, 102.32841856963, 0.0383001849303231, 0.04077752528139545),
        @("M1", "Ifnk", "Ifnar2", "Neutro", 1, 0.3333333333333333, 0.1379033333333333, 0.41371, 0.1151434304008604, 0.1151434304008603, 3, 1, 22.91393266666667, 68.741798, 0.09244477088342701, 0.09842430235492741, 3.15990769450889, 28.43916925058, 0.01064440804213936, 0.01133291180795782),
        @("M1", "Ifnk", "Ifnar2", "sCs", 1, 0.3333333333333333, 0.1379033333333333, 0.41371, 0.1151434304008604, 0.1151434304008603, 2, 1, 15.72803, 31.45606, 0.06345371398916565, 0.04503869334832845, 2.168947763766667, 13.0136865826, 0.007306278300387594, 0.005185909652898949),
        @("M2", "Ifnk", "Ifnar2", "ECs", 3, 1, 1.059762333333333, 3.179287, 0.8848565695991397, 0.8848565695991396, 2, 1, 29.4475095, 58.895019, 0.1188040616343076, 0.08432571340736816, 31.2073613785755, 187.244168271453, 0.1051245544321781, 0.07461616149464396),
        @("M2", "Ifnk", "Ifnar2", "FAPs", 3, 1, 1.059762333333333, 3.179287, 0.8848565695991397, 0.8848565695991396, 3, 1, 24.18563533333333, 72.556906, 0.09757537257289008, 0.1038867626663191, 25.63102533400244, 230.679228006022, 0.0863402094522055, 0.09192488443967906),
        @("M2", "Ifnk", "Ifnar2", "M1", 3, 1, 1.059762333333333, 3.179287, 0.8848565695991397, 0.8848565695991396, 3, 1, 73.143293, 219.429879, 0.2950918579280267, 0.3141790494976744, 77.51450685736366, 697.630561716273, 0.2611139691228304, 0.2780033959784305),
        @("M2", "Ifnk", "Ifnar2", "M2", 3, 1, 1.059762333333333, 3.179287, 0.8848565695991397, 0.8848565695991396, 3, 1, 82.44778433333333, 247.343353, 0.3326302229921831, 0.3541454787253825, 87.37505630325677, 786.3755067293109, 0.29433003806186, 0.313367953443987),
        @("M2", "Ifnk", "Ifnar2", "Neutro", 3, 1, 1.059762333333333, 3.179287, 0.8848565695991397, 0.8848565695991396, 3, 1, 22.91393266666667, 68.741798, 0.09244477088342701, 0.09842430235492741, 24.28332274866955, 218.549904738026, 0.08180036284128765, 0.08709139054696959),
        @("M2", "Ifnk", "Ifnar2", "sCs", 3, 1, 1.059762333333333, 3.179287, 0.8848565695991397, 0.8848565695991396, 2, 1, 15.72803, 31.45606, 0.06345371398916565, 0.04503869334832845, 16.66797377153667, 100.00784262922, 0.05614743568877806, 0.0398527836954295)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
